$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear cells whose value moved elsewhere in the row
$ws.Range("B13").ClearContents()
$ws.Range("C24").ClearContents()

# Cells that go from "X?" to "X" (done), or are brand new "X" (done) marks
# in column B (existing "X?" -> "X")
$ws.Range("B8").Value = "X"
$ws.Range("B30").Value = "X"
$ws.Range("B31").Value = "X"
$ws.Range("B34").Value = "X"
$ws.Range("B44").Value = "X"
$ws.Range("B50").Value = "X"
$ws.Range("B51").Value = "X"

# New "X" marks in column C
$ws.Range("C13").Value = "X"

# New "X" marks in column G (done column)
$ws.Range("G7").Value = "X"
$ws.Range("G8").Value = "X"
$ws.Range("G12").Value = "X"
$ws.Range("G13").Value = "X"
$ws.Range("G14").Value = "X"
$ws.Range("G24").Value = "X"
$ws.Range("G30").Value = "X"
$ws.Range("G31").Value = "X"
$ws.Range("G32").Value = "X"
$ws.Range("G34").Value = "X"
$ws.Range("G35").Value = "X"
$ws.Range("G36").Value = "X"
$ws.Range("G37").Value = "X"
$ws.Range("G38").Value = "X"
$ws.Range("G44").Value = "X"
$ws.Range("G50").Value = "X"
$ws.Range("G51").Value = "X"
$ws.Range("G53").Value = "X"

# New notes in column H
$ws.Range("H23").Value = "maybe wait for autoID branch"
$ws.Range("H48").Value = "maybe wait for autoID branch"

# Update selection to match the new active cell
$ws.Range("G15").Select()
